# Add a default header to the (single) section that shows the
# questionnaire number, e.g. "Questionnaire 39", centered, in Arial 12pt.

$d = $word.ActiveDocument

$section = $d.Sections.First

# wdHeaderFooterPrimary = 1  ->  the "default" header slot.
$header = $section.Headers.Item(1)

# Apply paragraph-level formatting (style + centered alignment) to the
# (currently empty) header range first.
$header.Range.Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Insert the header text. InsertAfter() (rather than assigning .Text)
# only materializes the single "default" header part/reference that we
# actually want, instead of also vivifying the even-page/first-page
# header & footer variants.
$header.Range.InsertAfter("Questionnaire 39")

# Re-grab the now-populated header range and trim off the trailing
# paragraph mark so the font formatting lands on the run only, not on
# the paragraph mark's run properties.
$full = $header.Range
$textRange = $full.Duplicate
$textRange.End = $full.End - 1

$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
